{"js": "// 1) Merge the two runs \"SHOOT FROM KEY AND \" + \"BACK UP TO BRIDGE\" in the\n//    \"AUTON SETTING 5\" heading paragraph into a single run with the same\n//    (bold/italic/underline) formatting.\nconst body = context.document.body;\n\nconst headingParas = body.paragraphs;\nheadingParas.load(\"items/text\");\nawait context.sync();\n\nconst mergedText = \"SHOOT FROM KEY AND BACK UP TO BRIDGE\";\nlet headingPara = null;\nfor (let i = 0; i < headingParas.items.length; i++) {\n  if (headingParas.items[i].text === mergedText) {\n    headingPara = headingParas.items[i];\n    break;\n  }\n}\nif (headingPara) {\n  // Re-writing the paragraph's text collapses it back down to a single run\n  // while keeping the run formatting (bold/italic/underline) already on it.\n  headingPara.insertText(mergedText, \"Replace\");\n  await context.sync();\n}\n\n// 2) Delete the whole \"6: Shoot from key at 2 point hoop\" paragraph.\nconst listParas = body.paragraphs;\nlistParas.load(\"items/text\");\nawait context.sync();\n\nlet targetPara = null;\nfor (let i = 0; i < listParas.items.length; i++) {\n  if (listParas.items[i].text.indexOf(\"Shoot from key at 2 point hoop\") >= 0) {\n    targetPara = listParas.items[i];\n    break;\n  }\n}\nif (targetPara) {\n  targetPara.delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Merge the two runs \"SHOOT FROM KEY AND \" + \"BACK UP TO BRIDGE\" (in the\n#    \"AUTON SETTING 5\" heading) into a single run. Word's Find/Replace\n#    engine matches across run boundaries and rewrites the hit as one run,\n#    carrying over the formatting (bold/italic/underline) already there.\n$mergedText = \"SHOOT FROM KEY AND BACK UP TO BRIDGE\"\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$findRange.Find.Replacement.ClearFormatting()\n$findRange.Find.Execute(\n    $mergedText,   # FindText\n    $true,         # MatchCase\n    $true,         # MatchWholeWord\n    $false,        # MatchWildcards\n    $false,        # MatchSoundsLike\n    $false,        # MatchAllWordForms\n    $true,         # Forward\n    1,             # Wrap (wdFindContinue)\n    $false,        # Format\n    $mergedText,   # ReplaceWith\n    2              # Replace (wdReplaceAll)\n) | Out-Null\n\n# 2) Delete the whole \"6: Shoot from key at 2 point hoop\" paragraph.\n$hitRange = $d.Content\n$found = $hitRange.Find.Execute(\"Shoot from key at 2 point hoop\")\nif ($found) {\n    $para = $hitRange.Paragraphs(1)\n    $para.Range.Delete()\n}\n"}
